$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phase 1: flush every old shared string used by B2:B51 out of the table by
# overwriting each cell with a unique throwaway value first. This guarantees
# none of the original strings remain referenced once we write the real data,
# so the shared-string table rebuilds cleanly in the exact order we write it.
$ws.Range("B2").Value = "__flush_2__"
$ws.Range("B3").Value = "__flush_3__"
$ws.Range("B4").Value = "__flush_4__"
$ws.Range("B5").Value = "__flush_5__"
$ws.Range("B6").Value = "__flush_6__"
$ws.Range("B7").Value = "__flush_7__"
$ws.Range("B8").Value = "__flush_8__"
$ws.Range("B9").Value = "__flush_9__"
$ws.Range("B10").Value = "__flush_10__"
$ws.Range("B11").Value = "__flush_11__"
$ws.Range("B12").Value = "__flush_12__"
$ws.Range("B13").Value = "__flush_13__"
$ws.Range("B14").Value = "__flush_14__"
$ws.Range("B15").Value = "__flush_15__"
$ws.Range("B16").Value = "__flush_16__"
$ws.Range("B17").Value = "__flush_17__"
$ws.Range("B18").Value = "__flush_18__"
$ws.Range("B19").Value = "__flush_19__"
$ws.Range("B20").Value = "__flush_20__"
$ws.Range("B21").Value = "__flush_21__"
$ws.Range("B22").Value = "__flush_22__"
$ws.Range("B23").Value = "__flush_23__"
$ws.Range("B24").Value = "__flush_24__"
$ws.Range("B25").Value = "__flush_25__"
$ws.Range("B26").Value = "__flush_26__"
$ws.Range("B27").Value = "__flush_27__"
$ws.Range("B28").Value = "__flush_28__"
$ws.Range("B29").Value = "__flush_29__"
$ws.Range("B30").Value = "__flush_30__"
$ws.Range("B31").Value = "__flush_31__"
$ws.Range("B32").Value = "__flush_32__"
$ws.Range("B33").Value = "__flush_33__"
$ws.Range("B34").Value = "__flush_34__"
$ws.Range("B35").Value = "__flush_35__"
$ws.Range("B36").Value = "__flush_36__"
$ws.Range("B37").Value = "__flush_37__"
$ws.Range("B38").Value = "__flush_38__"
$ws.Range("B39").Value = "__flush_39__"
$ws.Range("B40").Value = "__flush_40__"
$ws.Range("B41").Value = "__flush_41__"
$ws.Range("B42").Value = "__flush_42__"
$ws.Range("B43").Value = "__flush_43__"
$ws.Range("B44").Value = "__flush_44__"
$ws.Range("B45").Value = "__flush_45__"
$ws.Range("B46").Value = "__flush_46__"
$ws.Range("B47").Value = "__flush_47__"
$ws.Range("B48").Value = "__flush_48__"
$ws.Range("B49").Value = "__flush_49__"
$ws.Range("B50").Value = "__flush_50__"
$ws.Range("B51").Value = "__flush_51__"

# Phase 2: write the final values in the desired order.
$ws.Range("B2").Value = "_tejgfun_f5r08ct05pgrcopc"
$ws.Range("C2").Value = 0.01453988284457526
$ws.Range("B3").Value = "_tejgfun_f5r08ct05pgrco"
$ws.Range("C3").Value = 0.0113071124390496
$ws.Range("B4").Value = "_tejgfun_f2ct05amb"
$ws.Range("C4").Value = 0.01069627994486742
$ws.Range("B5").Value = "_tejgtotfun_f5r08pgrcopc"
$ws.Range("C5").Value = 0.01034545469611368
$ws.Range("B6").Value = "_tejgtotfun_f5r08pgrco"
$ws.Range("C6").Value = 0.00879744704735819
$ws.Range("B7").Value = "compu_muni_1"
$ws.Range("C7").Value = 0.006800414947068067
$ws.Range("B8").Value = "_tejgtotfun_f5amb"
$ws.Range("C8").Value = 0.00641436438821474
$ws.Range("B9").Value = "_tejgct_r08gstcrpc"
$ws.Range("C9").Value = 0.00625074298682951
$ws.Range("B10").Value = "_tejgge_r08ct05pobso"
$ws.Range("C10").Value = 0.005861190501612721
$ws.Range("B11").Value = "_tejgfun_f5ct05amb"
$ws.Range("C11").Value = 0.005454849117191334
$ws.Range("B12").Value = "_tejgrb_impmpc"
$ws.Range("C12").Value = 0.004733302525147709
$ws.Range("B13").Value = "_tejgrb_impm"
$ws.Range("C13").Value = 0.004477698035055916
$ws.Range("B14").Value = "_tejgge_r08ct05biser"
$ws.Range("C14").Value = 0.004010693001527983
$ws.Range("B15").Value = "_tejgct_r08gstcr"
$ws.Range("C15").Value = 0.003660237225034375
$ws.Range("B16").Value = "_tejgtotfun_f2opseg"
$ws.Range("C16").Value = 0.003476998729015072
$ws.Range("B17").Value = "_tejgfun_f2ct05opseg"
$ws.Range("C17").Value = 0.003362784643537152
$ws.Range("B18").Value = "_tejgtotfun_f2amb"
$ws.Range("C18").Value = 0.003311965964916121
$ws.Range("B19").Value = "_tejgct_r09gstcr"
$ws.Range("C19").Value = 0.002947250051367595
$ws.Range("B20").Value = "_tejgfun_f5ct05opseg"
$ws.Range("C20").Value = 0.002277719361433342
$ws.Range("B21").Value = "_tejgrb_redr"
$ws.Range("C21").Value = 0.002079869267472687
$ws.Range("B22").Value = "_tejgft_redr"
$ws.Range("C22").Value = 0.001813042162235562
$ws.Range("B23").Value = "_tejgtotfun_f2pgrco"
$ws.Range("C23").Value = 0.001661774444638138
$ws.Range("B24").Value = "_tejgge_r09ct05biser"
$ws.Range("C24").Value = 0.001607645658770075
$ws.Range("B25").Value = "_tejgfun_f5ct05prots"
$ws.Range("C25").Value = 0.001580438106153337
$ws.Range("B26").Value = "_tejgct_r09gstcrpc"
$ws.Range("C26").Value = 0.001434076594039504
$ws.Range("B27").Value = "tejgfun_f1ct05pgrco"
$ws.Range("C27").Value = 0.001410891486575873
$ws.Range("B28").Value = "devppimfun_f2ct06agro"
$ws.Range("C28").Value = 0.00138654968674241
$ws.Range("B29").Value = "_tejgfun_f2ct05pgrco"
$ws.Range("C29").Value = 0.001297504592780658
$ws.Range("B30").Value = "dfgpimpiafun_f2ct05agropc"
$ws.Range("C30").Value = 0.001266730292233024
$ws.Range("B31").Value = "tejgfun_f1ct05pgrcopc"
$ws.Range("C31").Value = 0.001212909696069222
$ws.Range("B32").Value = "piagtotfun_f5trans"
$ws.Range("C32").Value = 0.001210559054423836
$ws.Range("B33").Value = "pimgfun_f5r18ct06trans"
$ws.Range("C33").Value = 0.001152625747016905
$ws.Range("B34").Value = "devppimrb_canr"
$ws.Range("C34").Value = 0.001144557181629784
$ws.Range("B35").Value = "_tejgge_r08ct05biserpc"
$ws.Range("C35").Value = 0.001098145884837604
$ws.Range("B36").Value = "compu_muni_5"
$ws.Range("C36").Value = 0.00108690983935262
$ws.Range("B37").Value = "dfgpimpiafun_f5r08ct05agro"
$ws.Range("C37").Value = 0.001055598806744236
$ws.Range("B38").Value = "bvleche_16"
$ws.Range("C38").Value = 0.0009894476232293918
$ws.Range("B39").Value = "tdvgtotfun_f5r07comepc"
$ws.Range("C39").Value = 0.0009710239387325883
$ws.Range("B40").Value = "_devppimfun_f5ct06pgrco"
$ws.Range("C40").Value = 0.0009469991479901625
$ws.Range("B41").Value = "dfgdevpiagtotfun_f5r07comepc"
$ws.Range("C41").Value = 0.0009149755378173249
$ws.Range("B42").Value = "dfgdevpiagrb_foncpc"
$ws.Range("C42").Value = 0.0009117859812881071
$ws.Range("B43").Value = "pimgfun_f5r07ct05viv"
$ws.Range("C43").Value = 0.0009016316426015939
$ws.Range("B44").Value = "_tdvgfun_f5ct05prots"
$ws.Range("C44").Value = 0.0008622237936221261
$ws.Range("B45").Value = "_dfgdevpiagrb_foncpc"
$ws.Range("C45").Value = 0.0008368368683492897
$ws.Range("B46").Value = "devppimtotfun_f5r18opseg"
$ws.Range("C46").Value = 0.0008338826064467036
$ws.Range("B47").Value = "tejgtotfun_f1pgrco"
$ws.Range("C47").Value = 0.0008277056740861991
$ws.Range("B48").Value = "per_038"
$ws.Range("C48").Value = 0.0008274997601231334
$ws.Range("B49").Value = "dfgdevpiagge_r08ct05pobsopc"
$ws.Range("C49").Value = 0.0008224657154566134
$ws.Range("B50").Value = "_tejgge_r09ct05otgstpc"
$ws.Range("C50").Value = 0.0008208686560487466
$ws.Range("B51").Value = "piagtotfun_f5r07protspc"
$ws.Range("C51").Value = 0.0008073635729257993
